$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Promote the "flow_experiment_identifier" row (currently row 5) to become
#    the new row 2, right under the header -- shifting the old rows 2-4 down
#    by one (to rows 3-5). Row 6 onward is untouched.
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Copy()
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(6).Delete()

# Restore the row heights that belong with each row's wrapped text content
# (Insert() does not carry the source row height along).
$ws.Rows.Item(2).RowHeight = 102
$ws.Rows.Item(3).RowHeight = 136
$ws.Rows.Item(4).RowHeight = 170
$ws.Rows.Item(5).RowHeight = 170

# ---------------------------------------------------------------------------
# 2. Give the promoted row (row 2) a highlighted look: solid yellow fill
#    across the whole row.
# ---------------------------------------------------------------------------
$ws.Range("A2:R2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. Give the header row (row 1) a thin bottom border to separate it from
#    the data rows.
# ---------------------------------------------------------------------------
$hdr = $ws.Range("A1:R1")
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 4. Fix up hyperlinks: row 2's "Q2" hyperlink now belongs to "Q3" (where the
#    old row 2 content moved). Rows 6+ are unaffected. Rebuild the full
#    hyperlink collection since this runtime does not auto-shift hyperlink
#    anchors when rows are inserted/deleted.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("Q3"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q6:Q12"), "http://purl.obolibrary.org/obo/", "", "", "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("Q36:Q41"), "http://purl.obolibrary.org/obo/", "", "", "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("Q45:Q47"), "http://purl.obolibrary.org/obo/", "", "", "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("Q49"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q42"), "http://www.bioassayontology.org/bao/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q43"), "http://www.bioassayontology.org/bao/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q44"), "http://www.bioassayontology.org/bao/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q17:Q23"), "http://purl.obolibrary.org/obo/", "", "", "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("Q25:Q28"), "http://purl.obolibrary.org/obo/", "", "", "http://purl.obolibrary.org/obo/")
$ws.Hyperlinks.Add($ws.Range("Q13"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q14"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q15"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q16"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q26"), "http://purl.obolibrary.org/obo/", "", "", "")
$ws.Hyperlinks.Add($ws.Range("Q27"), "http://purl.obolibrary.org/obo/", "", "", "")

# ---------------------------------------------------------------------------
# 5. Reset the view: scroll back to the top-left and select a plain cell
#    instead of the whole column.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H56").Select()

Write-Output "done"
